$wb = $excel.ActiveWorkbook

$titles = @{
    "LTO7" = "LTO7 Internet Pricing as of 2022-08-08"
    "LTO8" = "LTO8 Internet Pricing as of 2022-08-08"
    "LTO9" = "LTO9 Internet Pricing as of 2022-08-08"
}

foreach ($ws in $wb.Worksheets) {
    $title = $titles[$ws.Name]
    if ($title) {
        $ws.Range("A1").Value = $title
    }
}
